$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# get_hotels now returns 3 hotels instead of 1 -- update the Rooms value
# for the sample/debug row to reflect a textual "1 Room" label.
$ws.Range("G2").Value = "1 Room"

# Debugging: move the active selection over to the Rooms column (G6).
$ws.Range("G6").Select()
